$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = 21470797556.071
$ws.Range("E4").Value  = 75490334188.022
$ws.Range("G4").Value  = 46767350124.022
$ws.Range("I4").Value  = 104213318252.022

$ws.Range("E7").Value  = 36829899326.955
$ws.Range("G7").Value  = 15280642796.445
$ws.Range("I7").Value  = 58379155857.465

$ws.Range("C12").Value = 7425556883.677
$ws.Range("E12").Value = 24009465280.814
$ws.Range("G12").Value = 14671539572.173
$ws.Range("I12").Value = 33347390989.455

$ws.Range("C13").Value = 9446555105.162001
